# Auto-generated edit script: updates market-price-derived leve profit columns (H-N)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching the
# upstream scheduled-runner market data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "ALC"; Row = 15; Cells = @{ "H" = 278.78946; "I" = 278.78946; "K" = 836.3683800000001; "M" = -667.3683800000001 } }
    @{ Sheet = "ALC"; Row = 69; Cells = @{ "H" = 20000; "J" = 20000; "L" = 60000; "N" = -61748 } }
    @{ Sheet = "ALC"; Row = 72; Cells = @{ "H" = 20000; "J" = 20000; "L" = 180000; "N" = -188736 } }
    @{ Sheet = "ALC"; Row = 107; Cells = @{ "H" = 790.25; "I" = 405; "J" = 918.6667; "K" = 405; "L" = 918.6667; "M" = 1515; "N" = -4758.6667 } }
    @{ Sheet = "ALC"; Row = 132; Cells = @{ "H" = 8999; "I" = 0; "K" = 0; "M" = $null } }
    @{ Sheet = "ALC"; Row = 137; Cells = @{ "H" = 914409.4; "I" = 1252187.5; "K" = 3756562.5; "M" = -3754012.5 } }
    @{ Sheet = "ARM"; Row = 2; Cells = @{ "H" = 0; "I" = 0; "K" = 0; "M" = $null } }
    @{ Sheet = "ARM"; Row = 5; Cells = @{ "H" = 300; "I" = 300; "J" = 0; "K" = 300; "L" = 0; "M" = -188; "N" = $null } }
    @{ Sheet = "ARM"; Row = 45; Cells = @{ "H" = 10820.714; "I" = 1915; "J" = 17500; "K" = 1915; "L" = 17500; "M" = -1538; "N" = -18254 } }
    @{ Sheet = "ARM"; Row = 63; Cells = @{ "H" = 7329.3335; "I" = 7000; "K" = 7000; "M" = -6314 } }
    @{ Sheet = "ARM"; Row = 66; Cells = @{ "H" = 7329.3335; "I" = 7000; "K" = 35000; "M" = -31568 } }
    @{ Sheet = "ARM"; Row = 69; Cells = @{ "H" = 179999; "J" = 179999; "L" = 179999; "N" = -181497 } }
    @{ Sheet = "ARM"; Row = 72; Cells = @{ "H" = 179999; "J" = 179999; "L" = 539997; "N" = -547485 } }
    @{ Sheet = "ARM"; Row = 116; Cells = @{ "H" = 0; "I" = 0; "K" = 0; "M" = $null } }
    @{ Sheet = "ARM"; Row = 131; Cells = @{ "H" = 50715; "J" = 50715; "L" = 50715; "N" = -60795 } }
    @{ Sheet = "BSM"; Row = 3; Cells = @{ "H" = 0; "I" = 0; "K" = 0; "M" = $null } }
    @{ Sheet = "BSM"; Row = 4; Cells = @{ "H" = 300; "I" = 300; "J" = 0; "K" = 300; "L" = 0; "M" = -185; "N" = $null } }
    @{ Sheet = "BSM"; Row = 22; Cells = @{ "H" = 250; "I" = 250; "K" = 250; "M" = -77 } }
    @{ Sheet = "BSM"; Row = 86; Cells = @{ "H" = 2400; "I" = 2400; "K" = 2400; "M" = -1277 } }
    @{ Sheet = "BSM"; Row = 89; Cells = @{ "H" = 2400; "I" = 2400; "K" = 12000; "M" = -6384 } }
    @{ Sheet = "BSM"; Row = 105; Cells = @{ "H" = 1771; "I" = 1688.75; "K" = 1688.75; "M" = 58.25 } }
    @{ Sheet = "BSM"; Row = 134; Cells = @{ "H" = 17749.75; "I" = 4000; "K" = 12000; "M" = -9465 } }
    @{ Sheet = "CRP"; Row = 5; Cells = @{ "H" = 243.6; "I" = 54.5; "J" = 1000; "K" = 54.5; "L" = 1000; "M" = 57.5; "N" = -1224 } }
    @{ Sheet = "CRP"; Row = 104; Cells = @{ "H" = 59999.5; "I" = 59999; "K" = 59999; "M" = -57378 } }
    @{ Sheet = "CRP"; Row = 107; Cells = @{ "H" = 556.61536; "I" = 614; "J" = 365.33334; "K" = 614; "L" = 365.33334; "M" = 1306; "N" = -4205.33334 } }
    @{ Sheet = "CUL"; Row = 17; Cells = @{ "H" = 67.666664; "J" = 1.5; "L" = 4.5; "N" = -342.5 } }
    @{ Sheet = "CUL"; Row = 23; Cells = @{ "H" = 150; "J" = 150; "L" = 450; "N" = -920 } }
    @{ Sheet = "CUL"; Row = 46; Cells = @{ "H" = 1937.5; "J" = 3375; "L" = 10125; "N" = -10307 } }
    @{ Sheet = "CUL"; Row = 81; Cells = @{ "H" = 403; "I" = 403; "K" = 1209; "M" = -86 } }
    @{ Sheet = "CUL"; Row = 84; Cells = @{ "H" = 403; "I" = 403; "K" = 3627; "M" = 1989 } }
    @{ Sheet = "CUL"; Row = 92; Cells = @{ "H" = 293.2; "I" = 293.2; "K" = 879.5999999999999; "M" = 368.4000000000001 } }
    @{ Sheet = "CUL"; Row = 113; Cells = @{ "H" = 314.85715; "I" = 250; "J" = 325.66666; "K" = 750; "L" = 976.9999799999999; "M" = 1420; "N" = -5316.99998 } }
    @{ Sheet = "CUL"; Row = 131; Cells = @{ "H" = 2303.4443; "J" = 2466.375; "L" = 7399.125; "N" = -17479.125 } }
    @{ Sheet = "GSM"; Row = 26; Cells = @{ "H" = 30500; "I" = 0; "J" = 30500; "K" = 0; "L" = 30500; "M" = $null; "N" = -31060 } }
    @{ Sheet = "GSM"; Row = 50; Cells = @{ "H" = 30500; "I" = 0; "J" = 30500; "K" = 0; "L" = 30500; "M" = $null; "N" = -31496 } }
    @{ Sheet = "GSM"; Row = 113; Cells = @{ "H" = 4583; "J" = 3999.5; "L" = 3999.5; "N" = -8339.5 } }
    @{ Sheet = "GSM"; Row = 132; Cells = @{ "H" = 123962.06; "I" = 169777.75; "K" = 509333.25; "M" = -506803.25 } }
    @{ Sheet = "LTW"; Row = 22; Cells = @{ "H" = 1333.3334; "I" = 1000; "J" = 1500; "K" = 1000; "L" = 1500; "M" = -705; "N" = -2090 } }
    @{ Sheet = "LTW"; Row = 27; Cells = @{ "H" = 1333.3334; "I" = 1000; "J" = 1500; "K" = 1000; "L" = 1500; "M" = -893; "N" = -1714 } }
    @{ Sheet = "LTW"; Row = 46; Cells = @{ "H" = 1000; "I" = 1000; "K" = 1000; "M" = -812 } }
    @{ Sheet = "LTW"; Row = 68; Cells = @{ "H" = 16001; "I" = 16001; "J" = 0; "K" = 16001; "L" = 0; "M" = -15252; "N" = $null } }
    @{ Sheet = "LTW"; Row = 71; Cells = @{ "H" = 16001; "I" = 16001; "J" = 0; "K" = 80005; "L" = 0; "M" = -76261; "N" = $null } }
    @{ Sheet = "LTW"; Row = 122; Cells = @{ "H" = 1999; "I" = 1748.75; "J" = 3000; "K" = 5246.25; "L" = 9000; "M" = -2796.25; "N" = -13900 } }
    @{ Sheet = "WVR"; Row = 2; Cells = @{ "H" = 501; "I" = 501; "J" = 0; "K" = 501; "L" = 0; "M" = -389; "N" = $null } }
    @{ Sheet = "WVR"; Row = 75; Cells = @{ "H" = 57500; "J" = 25000; "L" = 25000; "N" = -26872 } }
    @{ Sheet = "WVR"; Row = 78; Cells = @{ "H" = 57500; "J" = 25000; "L" = 75000; "N" = -84360 } }
    @{ Sheet = "WVR"; Row = 132; Cells = @{ "H" = 7320.4736; "I" = 6068.125; "K" = 18204.375; "M" = -15674.375 } }
    @{ Sheet = "WVR"; Row = 141; Cells = @{ "H" = 99856; "J" = 99856; "L" = 99856; "N" = -110216 } }
)

foreach ($entry in $updates) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($col in $entry.Cells.Keys) {
        $addr = "$col$($entry.Row)"
        $val = $entry.Cells[$col]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}

Write-Output "Applied $($updates.Count) row updates"
